$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawl re-run dropped 3 stale posts that used to sit at the top of the
# list (old rows 2-4). Remove them; every row below shifts up by three and
# the now-unused shared strings are dropped automatically on save.
$ws.Rows("2:4").Delete()

# The removed rows carried the cell that used to hold the hyperlink (old
# A6). That URL now lives in A3, so move the hyperlink there, restoring its
# original "visited" purple styling (the engine's own Hyperlinks.Add call
# applies a generic blue hyperlink font by default).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.facebook.com/congdongvnexpress/posts/919145590247453?ref=embed_post", "", "", "https://www.facebook.com/congdongvnexpress/posts/919145590247453?ref=embed_post")
$ws.Range("A3").Font.Color = 8388736

# Leave the selection where the editor last clicked after doing the cleanup.
$ws.Range("K13").Select()
